# ---------------------------------------------------------------------------
# Adds the 2022-Q3 quarterly snapshot:
#   1. A new worksheet "2022-Q3" (inserted right after "总计") with the
#      per-fund holdings table for that quarter.
#   2. A new row on the "总计" summary sheet for the 2022-Q3 totals, with the
#      existing quarters shifted down to make room (and renumbered).
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

function Set-HeaderCell($cell, $text) {
    # Header cells: bold, centered, thin border - matches the style used by
    # every other quarterly sheet's header row.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

function Set-IndexCell($cell, $num) {
    # Column-A row index cells share the header's bold/centered/bordered look.
    $cell.Value = $num
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

function Set-TextCell($cell, $text) {
    # Force text storage so numeric-looking values (fund codes with leading
    # zeros, percentages, decimals) are kept verbatim instead of becoming numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- Step 1: insert the new "2022-Q3" worksheet right after "总计" ------------
$sheetTotal = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $sheetTotal)
$newSheet.Name = "2022-Q3"

# Header row
Set-HeaderCell $newSheet.Cells.Item(1, 2) '基金代码'
Set-HeaderCell $newSheet.Cells.Item(1, 3) '基金名称'
Set-HeaderCell $newSheet.Cells.Item(1, 4) '基金规模'
Set-HeaderCell $newSheet.Cells.Item(1, 5) '股票总仓位'
Set-HeaderCell $newSheet.Cells.Item(1, 6) '仓位占比'
Set-HeaderCell $newSheet.Cells.Item(1, 7) '持有市值(亿元)'
Set-HeaderCell $newSheet.Cells.Item(1, 8) '仓位排名'

# Data rows (A: rank index, B-G: text fields incl. numeric-looking ones, H: numeric rank)
Set-IndexCell $newSheet.Cells.Item(2, 1) 0
Set-TextCell $newSheet.Cells.Item(2, 2) '012079'
Set-TextCell $newSheet.Cells.Item(2, 3) '信澳新能源精选混合'
Set-TextCell $newSheet.Cells.Item(2, 4) '50.03'
Set-TextCell $newSheet.Cells.Item(2, 5) '93.82'
Set-TextCell $newSheet.Cells.Item(2, 6) '4.61'
Set-TextCell $newSheet.Cells.Item(2, 7) '2.3064'
$newSheet.Cells.Item(2, 8).Value = 10
Set-IndexCell $newSheet.Cells.Item(3, 1) 1
Set-TextCell $newSheet.Cells.Item(3, 2) '007689'
Set-TextCell $newSheet.Cells.Item(3, 3) '国投瑞银新能源混合A'
Set-TextCell $newSheet.Cells.Item(3, 4) '43.07'
Set-TextCell $newSheet.Cells.Item(3, 5) '90.18'
Set-TextCell $newSheet.Cells.Item(3, 6) '4.79'
Set-TextCell $newSheet.Cells.Item(3, 7) '2.0631'
$newSheet.Cells.Item(3, 8).Value = 10
Set-IndexCell $newSheet.Cells.Item(4, 1) 2
Set-TextCell $newSheet.Cells.Item(4, 2) '007690'
Set-TextCell $newSheet.Cells.Item(4, 3) '国投瑞银新能源混合C'
Set-TextCell $newSheet.Cells.Item(4, 4) '36.94'
Set-TextCell $newSheet.Cells.Item(4, 5) '90.18'
Set-TextCell $newSheet.Cells.Item(4, 6) '4.79'
Set-TextCell $newSheet.Cells.Item(4, 7) '1.7694'
$newSheet.Cells.Item(4, 8).Value = 10
Set-IndexCell $newSheet.Cells.Item(5, 1) 3
Set-TextCell $newSheet.Cells.Item(5, 2) '005927'
Set-TextCell $newSheet.Cells.Item(5, 3) '创金合信新能源汽车主题股票A'
Set-TextCell $newSheet.Cells.Item(5, 4) '16.74'
Set-TextCell $newSheet.Cells.Item(5, 5) '92.78'
Set-TextCell $newSheet.Cells.Item(5, 6) '6.28'
Set-TextCell $newSheet.Cells.Item(5, 7) '1.0513'
$newSheet.Cells.Item(5, 8).Value = 9
Set-IndexCell $newSheet.Cells.Item(6, 1) 4
Set-TextCell $newSheet.Cells.Item(6, 2) '005928'
Set-TextCell $newSheet.Cells.Item(6, 3) '创金合信新能源汽车主题股票C'
Set-TextCell $newSheet.Cells.Item(6, 4) '16.65'
Set-TextCell $newSheet.Cells.Item(6, 5) '92.78'
Set-TextCell $newSheet.Cells.Item(6, 6) '6.28'
Set-TextCell $newSheet.Cells.Item(6, 7) '1.0456'
$newSheet.Cells.Item(6, 8).Value = 9
Set-IndexCell $newSheet.Cells.Item(7, 1) 5
Set-TextCell $newSheet.Cells.Item(7, 2) '005119'
Set-TextCell $newSheet.Cells.Item(7, 3) '银华智荟内在价值灵活配置混合A'
Set-TextCell $newSheet.Cells.Item(7, 4) '8.11'
Set-TextCell $newSheet.Cells.Item(7, 5) '93.91'
Set-TextCell $newSheet.Cells.Item(7, 6) '6.97'
Set-TextCell $newSheet.Cells.Item(7, 7) '0.5653'
$newSheet.Cells.Item(7, 8).Value = 4
Set-IndexCell $newSheet.Cells.Item(8, 1) 6
Set-TextCell $newSheet.Cells.Item(8, 2) '040001'
Set-TextCell $newSheet.Cells.Item(8, 3) '华安创新混合'
Set-TextCell $newSheet.Cells.Item(8, 4) '15.29'
Set-TextCell $newSheet.Cells.Item(8, 5) '72.21'
Set-TextCell $newSheet.Cells.Item(8, 6) '3.10'
Set-TextCell $newSheet.Cells.Item(8, 7) '0.4740'
$newSheet.Cells.Item(8, 8).Value = 9
Set-IndexCell $newSheet.Cells.Item(9, 1) 7
Set-TextCell $newSheet.Cells.Item(9, 2) '009859'
Set-TextCell $newSheet.Cells.Item(9, 3) '银华乐享混合A'
Set-TextCell $newSheet.Cells.Item(9, 4) '4.96'
Set-TextCell $newSheet.Cells.Item(9, 5) '94.27'
Set-TextCell $newSheet.Cells.Item(9, 6) '7.79'
Set-TextCell $newSheet.Cells.Item(9, 7) '0.3864'
$newSheet.Cells.Item(9, 8).Value = 5
Set-IndexCell $newSheet.Cells.Item(10, 1) 8
Set-TextCell $newSheet.Cells.Item(10, 2) '013160'
Set-TextCell $newSheet.Cells.Item(10, 3) '创金合信碳中和混合A'
Set-TextCell $newSheet.Cells.Item(10, 4) '5.11'
Set-TextCell $newSheet.Cells.Item(10, 5) '93.04'
Set-TextCell $newSheet.Cells.Item(10, 6) '7.37'
Set-TextCell $newSheet.Cells.Item(10, 7) '0.3766'
$newSheet.Cells.Item(10, 8).Value = 6
Set-IndexCell $newSheet.Cells.Item(11, 1) 9
Set-TextCell $newSheet.Cells.Item(11, 2) '001239'
Set-TextCell $newSheet.Cells.Item(11, 3) '长盛国企改革主题灵活配置混合'
Set-TextCell $newSheet.Cells.Item(11, 4) '4.46'
Set-TextCell $newSheet.Cells.Item(11, 5) '90.97'
Set-TextCell $newSheet.Cells.Item(11, 6) '7.81'
Set-TextCell $newSheet.Cells.Item(11, 7) '0.3483'
$newSheet.Cells.Item(11, 8).Value = 4
Set-IndexCell $newSheet.Cells.Item(12, 1) 10
Set-TextCell $newSheet.Cells.Item(12, 2) '005478'
Set-TextCell $newSheet.Cells.Item(12, 3) '长安鑫禧灵活配置混合C'
Set-TextCell $newSheet.Cells.Item(12, 4) '4.15'
Set-TextCell $newSheet.Cells.Item(12, 5) '91.79'
Set-TextCell $newSheet.Cells.Item(12, 6) '7.74'
Set-TextCell $newSheet.Cells.Item(12, 7) '0.3212'
$newSheet.Cells.Item(12, 8).Value = 9
Set-IndexCell $newSheet.Cells.Item(13, 1) 11
Set-TextCell $newSheet.Cells.Item(13, 2) '290002'
Set-TextCell $newSheet.Cells.Item(13, 3) '泰信先行策略混合'
Set-TextCell $newSheet.Cells.Item(13, 4) '6.17'
Set-TextCell $newSheet.Cells.Item(13, 5) '87.99'
Set-TextCell $newSheet.Cells.Item(13, 6) '5.13'
Set-TextCell $newSheet.Cells.Item(13, 7) '0.3165'
$newSheet.Cells.Item(13, 8).Value = 7
Set-IndexCell $newSheet.Cells.Item(14, 1) 12
Set-TextCell $newSheet.Cells.Item(14, 2) '013161'
Set-TextCell $newSheet.Cells.Item(14, 3) '创金合信碳中和混合C'
Set-TextCell $newSheet.Cells.Item(14, 4) '3.16'
Set-TextCell $newSheet.Cells.Item(14, 5) '93.04'
Set-TextCell $newSheet.Cells.Item(14, 6) '7.37'
Set-TextCell $newSheet.Cells.Item(14, 7) '0.2329'
$newSheet.Cells.Item(14, 8).Value = 6
Set-IndexCell $newSheet.Cells.Item(15, 1) 13
Set-TextCell $newSheet.Cells.Item(15, 2) '005076'
Set-TextCell $newSheet.Cells.Item(15, 3) '创金合信优选回报灵活配置混合'
Set-TextCell $newSheet.Cells.Item(15, 4) '2.74'
Set-TextCell $newSheet.Cells.Item(15, 5) '94.05'
Set-TextCell $newSheet.Cells.Item(15, 6) '8.42'
Set-TextCell $newSheet.Cells.Item(15, 7) '0.2307'
$newSheet.Cells.Item(15, 8).Value = 7
Set-IndexCell $newSheet.Cells.Item(16, 1) 14
Set-TextCell $newSheet.Cells.Item(16, 2) '004925'
Set-TextCell $newSheet.Cells.Item(16, 3) '长信低碳环保行业量化股票A'
Set-TextCell $newSheet.Cells.Item(16, 4) '6.08'
Set-TextCell $newSheet.Cells.Item(16, 5) '92.94'
Set-TextCell $newSheet.Cells.Item(16, 6) '3.41'
Set-TextCell $newSheet.Cells.Item(16, 7) '0.2073'
$newSheet.Cells.Item(16, 8).Value = 8
Set-IndexCell $newSheet.Cells.Item(17, 1) 15
Set-TextCell $newSheet.Cells.Item(17, 2) '121006'
Set-TextCell $newSheet.Cells.Item(17, 3) '国投瑞银稳健增长混合'
Set-TextCell $newSheet.Cells.Item(17, 4) '6.55'
Set-TextCell $newSheet.Cells.Item(17, 5) '61.70'
Set-TextCell $newSheet.Cells.Item(17, 6) '2.21'
Set-TextCell $newSheet.Cells.Item(17, 7) '0.1448'
$newSheet.Cells.Item(17, 8).Value = 10
Set-IndexCell $newSheet.Cells.Item(18, 1) 16
Set-TextCell $newSheet.Cells.Item(18, 2) '013842'
Set-TextCell $newSheet.Cells.Item(18, 3) '银华新锐成长混合A'
Set-TextCell $newSheet.Cells.Item(18, 4) '1.91'
Set-TextCell $newSheet.Cells.Item(18, 5) '94.78'
Set-TextCell $newSheet.Cells.Item(18, 6) '6.99'
Set-TextCell $newSheet.Cells.Item(18, 7) '0.1335'
$newSheet.Cells.Item(18, 8).Value = 3
Set-IndexCell $newSheet.Cells.Item(19, 1) 17
Set-TextCell $newSheet.Cells.Item(19, 2) '011147'
Set-TextCell $newSheet.Cells.Item(19, 3) '创金合信气候变化责任投资股票C'
Set-TextCell $newSheet.Cells.Item(19, 4) '1.43'
Set-TextCell $newSheet.Cells.Item(19, 5) '92.14'
Set-TextCell $newSheet.Cells.Item(19, 6) '8.74'
Set-TextCell $newSheet.Cells.Item(19, 7) '0.1250'
$newSheet.Cells.Item(19, 8).Value = 7
Set-IndexCell $newSheet.Cells.Item(20, 1) 18
Set-TextCell $newSheet.Cells.Item(20, 2) '015687'
Set-TextCell $newSheet.Cells.Item(20, 3) '银华乐享混合C'
Set-TextCell $newSheet.Cells.Item(20, 4) '1.39'
Set-TextCell $newSheet.Cells.Item(20, 5) '94.27'
Set-TextCell $newSheet.Cells.Item(20, 6) '7.79'
Set-TextCell $newSheet.Cells.Item(20, 7) '0.1083'
$newSheet.Cells.Item(20, 8).Value = 5
Set-IndexCell $newSheet.Cells.Item(21, 1) 19
Set-TextCell $newSheet.Cells.Item(21, 2) '011146'
Set-TextCell $newSheet.Cells.Item(21, 3) '创金合信气候变化责任投资股票A'
Set-TextCell $newSheet.Cells.Item(21, 4) '1.20'
Set-TextCell $newSheet.Cells.Item(21, 5) '92.14'
Set-TextCell $newSheet.Cells.Item(21, 6) '8.74'
Set-TextCell $newSheet.Cells.Item(21, 7) '0.1049'
$newSheet.Cells.Item(21, 8).Value = 7
Set-IndexCell $newSheet.Cells.Item(22, 1) 20
Set-TextCell $newSheet.Cells.Item(22, 2) '005477'
Set-TextCell $newSheet.Cells.Item(22, 3) '长安鑫禧灵活配置混合A'
Set-TextCell $newSheet.Cells.Item(22, 4) '1.23'
Set-TextCell $newSheet.Cells.Item(22, 5) '91.79'
Set-TextCell $newSheet.Cells.Item(22, 6) '7.74'
Set-TextCell $newSheet.Cells.Item(22, 7) '0.0952'
$newSheet.Cells.Item(22, 8).Value = 9
Set-IndexCell $newSheet.Cells.Item(23, 1) 21
Set-TextCell $newSheet.Cells.Item(23, 2) '290008'
Set-TextCell $newSheet.Cells.Item(23, 3) '泰信发展主题混合'
Set-TextCell $newSheet.Cells.Item(23, 4) '1.34'
Set-TextCell $newSheet.Cells.Item(23, 5) '91.43'
Set-TextCell $newSheet.Cells.Item(23, 6) '6.91'
Set-TextCell $newSheet.Cells.Item(23, 7) '0.0926'
$newSheet.Cells.Item(23, 8).Value = 5
Set-IndexCell $newSheet.Cells.Item(24, 1) 22
Set-TextCell $newSheet.Cells.Item(24, 2) '013104'
Set-TextCell $newSheet.Cells.Item(24, 3) '博时新能源主题混合C'
Set-TextCell $newSheet.Cells.Item(24, 4) '2.27'
Set-TextCell $newSheet.Cells.Item(24, 5) '85.72'
Set-TextCell $newSheet.Cells.Item(24, 6) '2.49'
Set-TextCell $newSheet.Cells.Item(24, 7) '0.0565'
$newSheet.Cells.Item(24, 8).Value = 9
Set-IndexCell $newSheet.Cells.Item(25, 1) 23
Set-TextCell $newSheet.Cells.Item(25, 2) '013103'
Set-TextCell $newSheet.Cells.Item(25, 3) '博时新能源主题混合A'
Set-TextCell $newSheet.Cells.Item(25, 4) '2.10'
Set-TextCell $newSheet.Cells.Item(25, 5) '85.72'
Set-TextCell $newSheet.Cells.Item(25, 6) '2.49'
Set-TextCell $newSheet.Cells.Item(25, 7) '0.0523'
$newSheet.Cells.Item(25, 8).Value = 9
Set-IndexCell $newSheet.Cells.Item(26, 1) 24
Set-TextCell $newSheet.Cells.Item(26, 2) '013843'
Set-TextCell $newSheet.Cells.Item(26, 3) '银华新锐成长混合C'
Set-TextCell $newSheet.Cells.Item(26, 4) '0.72'
Set-TextCell $newSheet.Cells.Item(26, 5) '94.78'
Set-TextCell $newSheet.Cells.Item(26, 6) '6.99'
Set-TextCell $newSheet.Cells.Item(26, 7) '0.0503'
$newSheet.Cells.Item(26, 8).Value = 3
Set-IndexCell $newSheet.Cells.Item(27, 1) 25
Set-TextCell $newSheet.Cells.Item(27, 2) '011273'
Set-TextCell $newSheet.Cells.Item(27, 3) '泰信景气驱动12个月持有期混合A'
Set-TextCell $newSheet.Cells.Item(27, 4) '1.09'
Set-TextCell $newSheet.Cells.Item(27, 5) '74.12'
Set-TextCell $newSheet.Cells.Item(27, 6) '4.57'
Set-TextCell $newSheet.Cells.Item(27, 7) '0.0498'
$newSheet.Cells.Item(27, 8).Value = 8
Set-IndexCell $newSheet.Cells.Item(28, 1) 26
Set-TextCell $newSheet.Cells.Item(28, 2) '013757'
Set-TextCell $newSheet.Cells.Item(28, 3) '泰信均衡价值混合A'
Set-TextCell $newSheet.Cells.Item(28, 4) '0.76'
Set-TextCell $newSheet.Cells.Item(28, 5) '66.30'
Set-TextCell $newSheet.Cells.Item(28, 6) '5.55'
Set-TextCell $newSheet.Cells.Item(28, 7) '0.0422'
$newSheet.Cells.Item(28, 8).Value = 3
Set-IndexCell $newSheet.Cells.Item(29, 1) 27
Set-TextCell $newSheet.Cells.Item(29, 2) '013151'
Set-TextCell $newSheet.Cells.Item(29, 3) '长信低碳环保行业量化股票C'
Set-TextCell $newSheet.Cells.Item(29, 4) '0.86'
Set-TextCell $newSheet.Cells.Item(29, 5) '92.94'
Set-TextCell $newSheet.Cells.Item(29, 6) '3.41'
Set-TextCell $newSheet.Cells.Item(29, 7) '0.0293'
$newSheet.Cells.Item(29, 8).Value = 8
Set-IndexCell $newSheet.Cells.Item(30, 1) 28
Set-TextCell $newSheet.Cells.Item(30, 2) '011274'
Set-TextCell $newSheet.Cells.Item(30, 3) '泰信景气驱动12个月持有期混合C'
Set-TextCell $newSheet.Cells.Item(30, 4) '0.42'
Set-TextCell $newSheet.Cells.Item(30, 5) '74.12'
Set-TextCell $newSheet.Cells.Item(30, 6) '4.57'
Set-TextCell $newSheet.Cells.Item(30, 7) '0.0192'
$newSheet.Cells.Item(30, 8).Value = 8
Set-IndexCell $newSheet.Cells.Item(31, 1) 29
Set-TextCell $newSheet.Cells.Item(31, 2) '013758'
Set-TextCell $newSheet.Cells.Item(31, 3) '泰信均衡价值混合C'
Set-TextCell $newSheet.Cells.Item(31, 4) '0.30'
Set-TextCell $newSheet.Cells.Item(31, 5) '66.30'
Set-TextCell $newSheet.Cells.Item(31, 6) '5.55'
Set-TextCell $newSheet.Cells.Item(31, 7) '0.0166'
$newSheet.Cells.Item(31, 8).Value = 3
Set-IndexCell $newSheet.Cells.Item(32, 1) 30
Set-TextCell $newSheet.Cells.Item(32, 2) '004244'
Set-TextCell $newSheet.Cells.Item(32, 3) '东方周期优选灵活配置混合'
Set-TextCell $newSheet.Cells.Item(32, 4) '0.35'
Set-TextCell $newSheet.Cells.Item(32, 5) '90.61'
Set-TextCell $newSheet.Cells.Item(32, 6) '4.50'
Set-TextCell $newSheet.Cells.Item(32, 7) '0.0158'
$newSheet.Cells.Item(32, 8).Value = 7
Set-IndexCell $newSheet.Cells.Item(33, 1) 31
Set-TextCell $newSheet.Cells.Item(33, 2) '016262'
Set-TextCell $newSheet.Cells.Item(33, 3) '银华智荟内在价值灵活配置混合C'
Set-TextCell $newSheet.Cells.Item(33, 4) '0.19'
Set-TextCell $newSheet.Cells.Item(33, 5) '93.91'
Set-TextCell $newSheet.Cells.Item(33, 6) '6.97'
Set-TextCell $newSheet.Cells.Item(33, 7) '0.0132'
$newSheet.Cells.Item(33, 8).Value = 4
Set-IndexCell $newSheet.Cells.Item(34, 1) 32
Set-TextCell $newSheet.Cells.Item(34, 2) '004360'
Set-TextCell $newSheet.Cells.Item(34, 3) '创金合信量化核心混合C'
Set-TextCell $newSheet.Cells.Item(34, 4) '0.21'
Set-TextCell $newSheet.Cells.Item(34, 5) '91.03'
Set-TextCell $newSheet.Cells.Item(34, 6) '1.42'
Set-TextCell $newSheet.Cells.Item(34, 7) '0.0030'
$newSheet.Cells.Item(34, 8).Value = 7
Set-IndexCell $newSheet.Cells.Item(35, 1) 33
Set-TextCell $newSheet.Cells.Item(35, 2) '004359'
Set-TextCell $newSheet.Cells.Item(35, 3) '创金合信量化核心混合A'
Set-TextCell $newSheet.Cells.Item(35, 4) '0.16'
Set-TextCell $newSheet.Cells.Item(35, 5) '91.03'
Set-TextCell $newSheet.Cells.Item(35, 6) '1.42'
Set-TextCell $newSheet.Cells.Item(35, 7) '0.0023'
$newSheet.Cells.Item(35, 8).Value = 7

# --- Step 2: update "总计" - insert a 2022-Q3 row, push older quarters down ----
$sheetTotal.Rows("2:2").Insert()
# The freshly inserted row inherits stray formatting from the Insert shift;
# clear it so the new row matches the plain look of the existing data rows.
$sheetTotal.Range("A2:D2").ClearFormats()

Set-IndexCell $sheetTotal.Cells.Item(2, 1) 0
$sheetTotal.Cells.Item(2, 2).Value = "2022-Q3"
$sheetTotal.Cells.Item(2, 3).Value = 34
$sheetTotal.Cells.Item(2, 4).Value = 12.85

# Renumber the index column (A) for the rows that shifted down so it stays 0-based sequential
for ($r = 3; $r -le 7; $r++) {
    $sheetTotal.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "2022-Q3 sheet added; 总计 summary updated."
